$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 13159.2
$ws.Range("I18").Value = 801
$ws.Range("J18").Value = 16248.75
$ws.Range("K18").Value = 801
$ws.Range("L18").Value = 16248.75
$ws.Range("M18").Value = -517
$ws.Range("N18").Value = -16816.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2216.8
$ws.Range("I40").Value = 2375
$ws.Range("J40").Value = 2111.3333
$ws.Range("K40").Value = 2375
$ws.Range("L40").Value = 2111.3333
$ws.Range("M40").Value = -2200
$ws.Range("N40").Value = -2461.3333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 529.75
$ws.Range("I96").Value = 529.75
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 1589.25
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -216.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 3293.889
$ws.Range("I98").Value = 3293.889
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 3293.889
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -1795.889

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 8879.267
$ws.Range("I113").Value = 10515.833
$ws.Range("J113").Value = 2333
$ws.Range("K113").Value = 10515.833
$ws.Range("L113").Value = 2333
$ws.Range("M113").Value = -7261.833000000001
$ws.Range("N113").Value = -8841

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 3293.889
$ws.Range("I122").Value = 3293.889
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 9881.667000000001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -7431.667000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 877.2558
$ws.Range("I129").Value = 880
$ws.Range("J129").Value = 876.8946999999999
$ws.Range("K129").Value = 2640
$ws.Range("L129").Value = 2630.6841
$ws.Range("M129").Value = 2360
$ws.Range("N129").Value = -12630.6841

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 969.9773
$ws.Range("I132").Value = 878.0714
$ws.Range("J132").Value = 2900
$ws.Range("K132").Value = 2634.2142
$ws.Range("L132").Value = 8700
$ws.Range("M132").Value = -104.2142000000003
$ws.Range("N132").Value = -13760

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1138.6666
$ws.Range("I2").Value = 1166.4
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 1166.4
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = -1053.4
$ws.Range("N2").Value = -1226

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1861.625
$ws.Range("I45").Value = 1600
$ws.Range("J45").Value = 1899
$ws.Range("K45").Value = 1600
$ws.Range("L45").Value = 1899
$ws.Range("M45").Value = -1223
$ws.Range("N45").Value = -2653

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 954
$ws.Range("I97").Value = 849.6667
$ws.Range("J97").Value = 1162.6666
$ws.Range("K97").Value = 849.6667
$ws.Range("L97").Value = 1162.6666
$ws.Range("M97").Value = -353.6667
$ws.Range("N97").Value = -2154.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1138.6666
$ws.Range("I116").Value = 1166.4
$ws.Range("J116").Value = 1000
$ws.Range("K116").Value = 1166.4
$ws.Range("L116").Value = 1000
$ws.Range("M116").Value = 1127.6
$ws.Range("N116").Value = -5588

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1138.6666
$ws.Range("I3").Value = 1166.4
$ws.Range("J3").Value = 1000
$ws.Range("K3").Value = 1166.4
$ws.Range("L3").Value = 1000
$ws.Range("M3").Value = -1052.4
$ws.Range("N3").Value = -1228

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 169558.83
$ws.Range("I86").Value = 2940.4
$ws.Range("J86").Value = 288572
$ws.Range("K86").Value = 2940.4
$ws.Range("L86").Value = 288572
$ws.Range("M86").Value = -1817.4
$ws.Range("N86").Value = -290818

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 169558.83
$ws.Range("I89").Value = 2940.4
$ws.Range("J89").Value = 288572
$ws.Range("K89").Value = 14702
$ws.Range("L89").Value = 1442860
$ws.Range("M89").Value = -9086
$ws.Range("N89").Value = -1454092

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 233.81818
$ws.Range("I7").Value = 173.75
$ws.Range("J7").Value = 394
$ws.Range("K7").Value = 173.75
$ws.Range("L7").Value = 394
$ws.Range("M7").Value = -60.75
$ws.Range("N7").Value = -620

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1249.8572
$ws.Range("I22").Value = 300
$ws.Range("J22").Value = 1408.1666
$ws.Range("K22").Value = 300
$ws.Range("L22").Value = 1408.1666
$ws.Range("M22").Value = 50
$ws.Range("N22").Value = -2108.1666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1634.4783
$ws.Range("I31").Value = 1198.9375
$ws.Range("J31").Value = 2630
$ws.Range("K31").Value = 1198.9375
$ws.Range("L31").Value = 2630
$ws.Range("M31").Value = -903.9375
$ws.Range("N31").Value = -3220

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1634.4783
$ws.Range("I34").Value = 1198.9375
$ws.Range("J34").Value = 2630
$ws.Range("K34").Value = 1198.9375
$ws.Range("L34").Value = 2630
$ws.Range("M34").Value = -996.9375
$ws.Range("N34").Value = -3034

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1591.0605
$ws.Range("I58").Value = 929.86957
$ws.Range("J58").Value = 3111.8
$ws.Range("K58").Value = 929.86957
$ws.Range("L58").Value = 3111.8
$ws.Range("M58").Value = -726.86957
$ws.Range("N58").Value = -3517.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2393.4666
$ws.Range("I99").Value = 1581.5
$ws.Range("J99").Value = 2934.7778
$ws.Range("K99").Value = 1581.5
$ws.Range("L99").Value = 2934.7778
$ws.Range("M99").Value = -83.5
$ws.Range("N99").Value = -5930.7778

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2393.4666
$ws.Range("I126").Value = 1581.5
$ws.Range("J126").Value = 2934.7778
$ws.Range("K126").Value = 4744.5
$ws.Range("L126").Value = 8804.3334
$ws.Range("M126").Value = -2274.5
$ws.Range("N126").Value = -13744.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1761.2812
$ws.Range("I134").Value = 1516.5
$ws.Range("J134").Value = 3474.75
$ws.Range("K134").Value = 4549.5
$ws.Range("L134").Value = 10424.25
$ws.Range("M134").Value = -2014.5
$ws.Range("N134").Value = -15494.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1591.0605
$ws.Range("I136").Value = 929.86957
$ws.Range("J136").Value = 3111.8
$ws.Range("K136").Value = 2789.60871
$ws.Range("L136").Value = 9335.400000000001
$ws.Range("M136").Value = -239.60871
$ws.Range("N136").Value = -14435.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 422.4
$ws.Range("I38").Value = 70
$ws.Range("J38").Value = 951
$ws.Range("K38").Value = 210
$ws.Range("L38").Value = 2853
$ws.Range("M38").Value = 137
$ws.Range("N38").Value = -3547

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1097.6666
$ws.Range("I122").Value = 899.6667
$ws.Range("J122").Value = 1196.6666
$ws.Range("K122").Value = 8097.0003
$ws.Range("L122").Value = 10769.9994
$ws.Range("M122").Value = -5647.0003
$ws.Range("N122").Value = -15669.9994

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 779.3608400000001
$ws.Range("I131").Value = 521
$ws.Range("J131").Value = 796.3955999999999
$ws.Range("K131").Value = 1563
$ws.Range("L131").Value = 2389.1868
$ws.Range("M131").Value = 3477
$ws.Range("N131").Value = -12469.1868

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1038.6
$ws.Range("I80").Value = 997.6667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 1038.6
$ws.Range("I83").Value = 997.6667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2603.4546
$ws.Range("I132").Value = 2189.8845
$ws.Range("J132").Value = 4139.5713
$ws.Range("K132").Value = 6569.6535
$ws.Range("L132").Value = 12418.7139
$ws.Range("M132").Value = -4039.6535
$ws.Range("N132").Value = -17478.7139

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4242.5293
$ws.Range("I40").Value = 1941.3
$ws.Range("J40").Value = 7530
$ws.Range("K40").Value = 1941.3
$ws.Range("L40").Value = 7530
$ws.Range("M40").Value = -1805.3
$ws.Range("N40").Value = -7802

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2760.4443
$ws.Range("I46").Value = 1400
$ws.Range("J46").Value = 3149.1428
$ws.Range("K46").Value = 1400
$ws.Range("L46").Value = 3149.1428
$ws.Range("M46").Value = -1212
$ws.Range("N46").Value = -3525.1428

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 3096.5
$ws.Range("I82").Value = 2498.3333
$ws.Range("J82").Value = 3993.75
$ws.Range("K82").Value = 2498.3333
$ws.Range("L82").Value = 3993.75
$ws.Range("M82").Value = -2137.3333
$ws.Range("N82").Value = -4715.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 3096.5
$ws.Range("I85").Value = 2498.3333
$ws.Range("J85").Value = 3993.75
$ws.Range("K85").Value = 2498.3333
$ws.Range("L85").Value = 3993.75
$ws.Range("M85").Value = -1250.3333
$ws.Range("N85").Value = -6489.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 763.1667
$ws.Range("I93").Value = 800
$ws.Range("J93").Value = 689.5
$ws.Range("K93").Value = 800
$ws.Range("L93").Value = 689.5
$ws.Range("M93").Value = 448
$ws.Range("N93").Value = -3185.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1801.4286
$ws.Range("I100").Value = 1601.6666
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 1601.6666
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -1060.6666
$ws.Range("N100").Value = -4082

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H53").Value = 7500
$ws.Range("I53").Value = 7500
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 7500
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = -6893

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2052.8572
$ws.Range("I81").Value = 2052.8572
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 4105.7144
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -3044.7144

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 2052.8572
$ws.Range("I84").Value = 2052.8572
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 20528.572
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -15224.572

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2866.6667
$ws.Range("I96").Value = 600
$ws.Range("J96").Value = 4000
$ws.Range("K96").Value = 600
$ws.Range("L96").Value = 4000
$ws.Range("M96").Value = 773
$ws.Range("N96").Value = -6746

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 34754.434
$ws.Range("I122").Value = 37979.19
$ws.Range("J122").Value = 894.5
$ws.Range("K122").Value = 113937.57
$ws.Range("L122").Value = 2683.5
$ws.Range("M122").Value = -111487.57
$ws.Range("N122").Value = -7583.5
